$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) -------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value  = 55
$ws1.Range("F4").Value  = 16323
$ws1.Range("G4").Value  = 70
$ws1.Range("G5").Value  = "不可售"
$ws1.Range("F6").Value  = 19
$ws1.Range("F7").Value  = 733
$ws1.Range("F8").Value  = 15592
$ws1.Range("G10").Value = 80
$ws1.Range("F11").Value = 463
$ws1.Range("F13").Value = 1026
$ws1.Range("F17").Value = 221
$ws1.Range("F19").Value = 89
$ws1.Range("F20").Value = 606
$ws1.Range("F28").Value = 523
$ws1.Range("F29").Value = 38
$ws1.Range("F32").Value = 81
$ws1.Range("F36").Value = 363
$ws1.Range("F39").Value = 5674
$ws1.Range("F40").Value = 5243

# --- Sheet "演出" (Performance) ------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3

# --- Sheet "全部类型" (All types) ----------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value  = 55
$ws4.Range("F4").Value  = 16323
$ws4.Range("G4").Value  = 70
$ws4.Range("G5").Value  = "不可售"
$ws4.Range("F6").Value  = 19
$ws4.Range("F7").Value  = 733
$ws4.Range("F8").Value  = 15592
$ws4.Range("G10").Value = 80
$ws4.Range("F11").Value = 463
$ws4.Range("F13").Value = 1026
$ws4.Range("F17").Value = 221
$ws4.Range("F19").Value = 89
$ws4.Range("F20").Value = 606
$ws4.Range("F28").Value = 523
$ws4.Range("F29").Value = 38
$ws4.Range("F34").Value = 81
$ws4.Range("F38").Value = 363
$ws4.Range("F41").Value = 5674
$ws4.Range("F42").Value = 3
$ws4.Range("F43").Value = 5243
